$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Test_Run_YN value for the Login_To_App scenario (row 2) from "N" to "Y"
$ws.Range("D2").Value = "Y"

# Move the active selection to D3, matching the post-edit sheet view
$ws.Range("D3").Select()
